# biter coin values.xlsx — commit: "added biter killing bounty / fixed science
# recipe missing names and animation oddity"
#
# The "l" and "b" biter-size codes in column A were placeholders; give them
# their real Factorio names. Also re-enter the D column (square-root) formula
# as one fill so it becomes a single shared formula across D2:D8 (matches how
# the other derived columns, C/E/F/G/H, are already filled) and nudge the
# view state (column H width + selected cell) to match the saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- fixed science recipe missing names ---------------------------------
$ws.Range("A4").Value = "big"
$ws.Range("A5").Value = "behemoth"

# --- fixed "animation oddity": D2:D8 re-filled as one shared formula ----
$ws.Range("D2:D8").Formula = "=SQRT(B2)/2"

# --- cosmetic view-state touch-up ----------------------------------------
$ws.Columns.Item(8).ColumnWidth = 12.67
$ws.Range("F18").Select()
